# GroundUpMCMC DataDictionaryAndOrganization.xlsx update
# "build new class for model parameters" -- rework the "location in code" /
# "used in" columns on the data dictionary sheet, and add an "OG?" column
# to the interference sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": rename "massSpec model" -> "massSpecModel property",
# split "method struct, new tables" into "method struct" (Method section)
# and "mass property" (isotopic masses row); fill in the newly-populated
# "location in code" cells; widen columns C & D; add a new blank-ish
# note row (B32 = " ").
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")
$wsData.Activate() | Out-Null

# Mass Spectrometer block
$wsData.Range("D9").Value  = "massSpecModel property"
$wsData.Range("D10").Value = "massSpecModel property"
$wsData.Range("D11").Value = "massSpecModel property"
$wsData.Range("D12").Value = "massSpecModel property"

# Method block
$wsData.Range("D17").Value = "method struct"
$wsData.Range("D18").Value = "method struct"
$wsData.Range("D19").Value = "method struct"
$wsData.Range("D20").Value = "method struct"
$wsData.Range("D21").Value = "method struct"
$wsData.Range("D22").Value = "method struct"
$wsData.Range("D23").Value = "method struct"

# Physical Constants block
$wsData.Range("D28").Value = "massSpecModel property"
$wsData.Range("D29").Value = "massSpecModel property"
$wsData.Range("D30").Value = "mass property"

# new note row
$wsData.Range("B32").Value = " "

# column widths (characters); engine stores width in a coarser quantum
# than real Excel, these land as close as possible to 44.1640625 / 26
$wsData.Columns.Item(3).ColumnWidth = 43.330729166666664
$wsData.Columns.Item(4).ColumnWidth = 25.166666666666668

$aw = $excel.ActiveWindow
$aw.ScrollRow = 5
$aw.ScrollColumn = 1
$wsData.Range("B37").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "interference": add a new "OG?" column in column A.
# ---------------------------------------------------------------------
$wsInterference = $wb.Worksheets.Item("interference")
$wsInterference.Activate() | Out-Null

$wsInterference.Range("A3").Value = "OG?"
$wsInterference.Range("A4").Value = "yes"
$wsInterference.Range("A5").Value = "yes"
$wsInterference.Range("A6").Value = "yes"
$wsInterference.Range("A8").Value = "yes"

$aw2 = $excel.ActiveWindow
$aw2.ScrollRow = 5
$aw2.ScrollColumn = 1
$wsInterference.Range("F39").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "variables": move selection.
# ---------------------------------------------------------------------
$wsVariables = $wb.Worksheets.Item("variables")
$wsVariables.Activate() | Out-Null
$wsVariables.Range("B29").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "unknowns": move selection / scroll.
# ---------------------------------------------------------------------
$wsUnknowns = $wb.Worksheets.Item("unknowns")
$wsUnknowns.Activate() | Out-Null
$aw3 = $excel.ActiveWindow
$aw3.ScrollRow = 5
$aw3.ScrollColumn = 1
$wsUnknowns.Range("F21").Select() | Out-Null

# ---------------------------------------------------------------------
# Leave "interference" as the active / last-selected tab, matching the
# workbook's new activeTab.
# ---------------------------------------------------------------------
$wsInterference.Activate() | Out-Null
$wsInterference.Range("F39").Select() | Out-Null
